# RegistrationTestData.xlsx edit:
#   The test-data email address used for the "EmailAddress" and
#   "ConfirmEmailAddress" columns on the "Registration" sheet was bumped
#   from ravitest64@yopmail.com to ravitest67@yopmail.com.
#
# Both cells (I2 and J2) share the same underlying string value, so both
# are updated to keep the registration test data internally consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registration")

$oldEmail = "ravitest64@yopmail.com"
$newEmail = "ravitest67@yopmail.com"

$emailCell = $ws.Range("I2")
if ($emailCell.Value2 -eq $oldEmail) {
    $emailCell.Value2 = $newEmail
}

$confirmEmailCell = $ws.Range("J2")
if ($confirmEmailCell.Value2 -eq $oldEmail) {
    $confirmEmailCell.Value2 = $newEmail
}
